# Auto-generated edit script: updates Price (D) and Volume(1h) (E) columns
# on the cryptos sheet to reflect the refreshed scrape.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '66.875.81'
$ws.Range('E2').Value = '  +0.00%  '
$ws.Range('D3').Value = '3.112.00'
$ws.Range('E3').Value = '  +0.61%  '
$ws.Range('E4').Value = '  +0.02%  '
# D5: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '577.27'
$ws.Range('D5').ClearFormats()
# D6: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '172.94'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  +2.61%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('E8').Value = '  -0.42%  '
$ws.Range('E9').Value = '  -2.99%  '
$ws.Range('E10').Value = '  -1.70%  '
$ws.Range('E11').Value = '  -0.07%  '
$ws.Range('E12').Value = '  -1.08%  '
# D13: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '37.15'
$ws.Range('D13').ClearFormats()
$ws.Range('E13').Value = '  +2.00%  '
$ws.Range('E14').Value = '  -1.33%  '
$ws.Range('D15').Value = '3.630.19'
$ws.Range('E15').Value = '  +0.70%  '
$ws.Range('D16').Value = '66.847.16'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('D18').Value = '3.114.54'
$ws.Range('E18').Value = '  +0.70%  '
# D19: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '16.20'
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.37%  '
# D20: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '475.26'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.72%  '
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +4.76%  '
# D23: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '83.89'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.10%  '
# D24: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.24'
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  +0.87%  '
$ws.Range('E25').Value = '  -3.74%  '
# D26: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '10.11'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  -0.41%  '
# D28: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '7.94'
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  -0.60%  '
$ws.Range('E29').Value = '  -1.16%  '
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('E31').Value = '  +1.00%  '
# D32: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.115'
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  -0.10%  '
$ws.Range('D33').Value = '0.0₃0952'
$ws.Range('E33').Value = '  -7.42%  '
# D34: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  +0.01%  '
$ws.Range('E35').Value = '  -1.15%  '
# D36: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.974'
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.10%  '
# D37: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '47.10'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +0.25%  '
# D38: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '50.18'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  -0.26%  '
$ws.Range('E39').Value = '  -2.40%  '
$ws.Range('E40').Value = '  -2.27%  '
# D41: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.123'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +1.68%  '
# D42: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '8.66'
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.24%  '
$ws.Range('D43').Value = '2.811.89'
$ws.Range('E43').Value = '  +1.18%  '
# D44: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '381.72'
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.35%  '
$ws.Range('E45').Value = '  -2.09%  '
$ws.Range('E46').Value = '  -9.94%  '
# D47: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '135.41'
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  +0.12%  '
# D49: numeric-looking text -> force text type, then restore default (no) style
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '24.79'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -0.90%  '
$ws.Range('E50').Value = '  -1.54%  '
$ws.Range('E51').Value = '  -0.83%  '
